# Append " -- DONE" / " --DONE" as a brand-new trailing run to the eight
# To-Do paragraphs that don't yet have it, mirroring the two paragraphs
# that already carry a "--DONE" suffix.
#
# We work on the whole document body's OOXML at once (rather than
# per-paragraph Range.InsertXML) because the very last paragraph in the
# body has no paragraph-mark character inside its own Range here, which
# makes per-paragraph XML splicing unreliable right at the end of the
# document. Pulling/pushing the WordOpenXML for $d.Content sidesteps that.

$d = $word.ActiveDocument

$full = $d.Content.WordOpenXML
$m = [regex]::Match(
    $full,
    '(?s)<pkg:part pkg:name="/word/document\.xml"[^>]*>.*?<pkg:xmlData>(?<doc>.*?)</pkg:xmlData>\s*</pkg:part>'
)
if (-not $m.Success) { throw "could not locate document.xml part in WordOpenXML" }
$docXml = $m.Groups["doc"].Value

# $d.Content.WordOpenXML tacks on a placeholder empty paragraph right
# before <w:sectPr> as an artifact of serializing the whole-document
# range; drop it so we don't introduce a spurious extra paragraph.
$docXml = $docXml -replace '<w:p w14:paraId="00000001" w14:textId="77777777" w:rsidR="00000000" w:rsidRDefault="00000000"/>', ''

function Add-DoneRun([string]$xml, [string]$marker, [string]$suffix) {
    $idx = $xml.IndexOf($marker)
    if ($idx -lt 0) { throw "marker not found: $marker" }
    $runClose = '</w:t></w:r>'
    $closeIdx = $xml.IndexOf($runClose, $idx)
    if ($closeIdx -lt 0) { throw "run close not found after marker: $marker" }
    $insertPos = $closeIdx + $runClose.Length
    $newRun = '<w:r><w:t xml:space="preserve">' + $suffix + '</w:t></w:r>'
    return $xml.Substring(0, $insertPos) + $newRun + $xml.Substring($insertPos)
}

$docXml = Add-DoneRun $docXml 'Add labels to the chunk that computes CDRs' ' -- DONE'
$docXml = Add-DoneRun $docXml 'Make comments more expressive.' ' -- DONE'
$docXml = Add-DoneRun $docXml 'Add comments to the last chunk to explain what it is doing.' ' -- DONE'
$docXml = Add-DoneRun $docXml 'Fix the last comment sentence.' ' --DONE'
$docXml = Add-DoneRun $docXml 'Add text and comments in the chunk compute-change' ' -- DONE'
$docXml = Add-DoneRun $docXml 'Get rid of the t-test' ' -- DONE'
$docXml = Add-DoneRun $docXml 'When calculating ASFRs, include appropriate age range only.' ' --DONE'
$docXml = Add-DoneRun $docXml 'Flesh out responses' ' --DONE'

$wrapped = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' + $docXml + '</pkg:xmlData></pkg:part></pkg:package>'

$d.Content.InsertXML($wrapped) | Out-Null

Write-Output "Appended DONE markers to the To Do list."
